# "Spring Hibernate Integration in backend project"
# Fill in row 5 (sheet index 5, i.e. S.No = 5) of the DAY_01 log sheet with the
# new Spring Hibernate task, and update the sheet's selection / new helper
# column that Excel recorded as a side effect of the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # DAY_01 - the active/selected sheet in the workbook

# --- Row 6 (S.No 5): Spring Hibernate configuration task -------------------
$ws.Range("B6").Value = "17/102016"
$ws.Range("C6").Value = "Spring Hibernate configuration"
$ws.Range("E6").Value = "yes"
$ws.Range("F6").Value = "2hr"
$ws.Range("G6").Value = "Y"

$fullError = "INFO: Using DataSource [org.springframework.jdbc.datasource.DriverManagerDataSource@31d7b7bf] of Hibernate SessionFactory for HibernateTransactionManager`nException in thread ""main"" javax.persistence.PersistenceException: org.hibernate.PersistentObjectException: detached entity passed to persist: com.niit.Backend.Model.ProductBean"
$ws.Range("H6").Value = $fullError
$ws.Range("H6").WrapText = $true

# --- New (empty) helper column I, sized like the rest of the table ---------
$ws.Columns.Item(9).ColumnWidth = 14.33

# --- Selection moved to C7 (top-left cell scroll reset) ---------------------
$ws.Range("C7").Select()
